$d = $word.ActiveDocument

# --- Paragraph "ΤΥ47 - Αξιολόγηση από τον Πελάτη" -> "... πελάτη" ---
$p47 = $d.Paragraphs(47)
$p47Start = $p47.Range.Start
$p47Text = $p47.Range.Text
$idx47Pi = $p47Text.IndexOf("Π")
$idx47Num = $p47Text.IndexOf("7")

$r47Pi = $d.Range($p47Start + $idx47Pi, $p47Start + $idx47Pi + 1)
$r47Pi.Text = "π"

$r47Num = $d.Range($p47Start + $idx47Num, $p47Start + $idx47Num + 1)
$r47Num.Font.Bold = $true
$r47Num.Font.Bold = $false

$r47PiAgain = $d.Range($p47Start + $idx47Pi, $p47Start + $idx47Pi + 1)
$r47PiAgain.Font.Bold = $true
$r47PiAgain.Font.Bold = $false

# --- Paragraph "ΤΥ48 – Τροποποιήσεις Συστήματος" -> "... συστήματος" ---
$p48 = $d.Paragraphs(48)
$p48Start = $p48.Range.Start
$p48Text = $p48.Range.Text
$idx48Sigma = $p48Text.IndexOf("Σ")
$idx48Num = $p48Text.IndexOf("8")

$r48Sigma = $d.Range($p48Start + $idx48Sigma, $p48Start + $idx48Sigma + 1)
$r48Sigma.Text = "σ"

$r48Num = $d.Range($p48Start + $idx48Num, $p48Start + $idx48Num + 1)
$r48Num.Font.Bold = $true
$r48Num.Font.Bold = $false

$r48SigmaAgain = $d.Range($p48Start + $idx48Sigma, $p48Start + $idx48Sigma + 1)
$r48SigmaAgain.Font.Bold = $true
$r48SigmaAgain.Font.Bold = $false

# --- Paragraph "ΤΥ49 – Αποδοχή και Εγκατάσταση" -> "... εγκατάσταση" with the
#     _GoBack bookmark moved here, splitting "ε" from "γκατάσταση" ---
$p49 = $d.Paragraphs(49)
$p49Start = $p49.Range.Start
$p49Text = $p49.Range.Text
$idx49Eps = $p49Text.IndexOf("Ε")
$idx49Num = $p49Text.IndexOf("49")

$r49Eps = $d.Range($p49Start + $idx49Eps, $p49Start + $idx49Eps + 1)
$r49Eps.Text = "ε"

$r49Num = $d.Range($p49Start + $idx49Num, $p49Start + $idx49Num + 2)
$r49Num.Font.Bold = $true
$r49Num.Font.Bold = $false

$r49EpsAgain = $d.Range($p49Start + $idx49Eps, $p49Start + $idx49Eps + 1)
$r49EpsAgain.Font.Bold = $true
$r49EpsAgain.Font.Bold = $false

# Point (zero-length) range right after the lower-cased epsilon -> add the
# _GoBack bookmark there. Since bookmark names are unique, this also removes
# the old _GoBack that used to sit at the end of the "...Δεδομένων" paragraph.
$bmPos = $p49Start + $idx49Eps + 1
$bmRange = $d.Range($bmPos, $bmPos)
$bmRange.Bookmarks.Add("_GoBack")

"done"
